$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2775
$ws.Range("J3").Value = 2871
$ws.Range("B4").Value = 1675
$ws.Range("D4").Value = 1952
$ws.Range("I4").Value = 1757
$ws.Range("J4").Value = 648
$ws.Range("J6").Value = 3509
$ws.Range("B7").Value = 23307
$ws.Range("D7").Value = 28142
$ws.Range("I7").Value = 26204
$ws.Range("J7").Value = 10023

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 101
$ws.Range("J3").Value = 122
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 47
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 76
$ws.Range("J7").Value = 363

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 78
$ws.Range("J5").Value = 26
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 301
$ws.Range("J8").Value = 635
$ws.Range("J9").Value = 62
$ws.Range("J10").Value = 62
$ws.Range("J11").Value = 144
$ws.Range("I18").Value = 206
$ws.Range("J19").Value = 317
$ws.Range("J20").Value = 207
$ws.Range("J21").Value = 17
$ws.Range("J23").Value = 103
$ws.Range("J25").Value = 57
$ws.Range("J29").Value = 576
$ws.Range("J31").Value = 75
$ws.Range("J33").Value = 414
$ws.Range("J34").Value = 53
$ws.Range("J37").Value = 338
$ws.Range("J42").Value = 392
$ws.Range("J46").Value = 35
$ws.Range("J48").Value = 101
$ws.Range("J50").Value = 57
$ws.Range("J51").Value = 134
$ws.Range("J52").Value = 255
$ws.Range("J53").Value = 98
$ws.Range("J54").Value = 198
$ws.Range("J60").Value = 66
$ws.Range("B63").Value = 379
$ws.Range("D63").Value = 334
$ws.Range("J63").Value = 47
$ws.Range("J67").Value = 363
$ws.Range("J73").Value = 91
$ws.Range("J75").Value = 31
$ws.Range("J77").Value = 85
$ws.Range("J78").Value = 134
$ws.Range("J79").Value = 300
$ws.Range("J80").Value = 20
$ws.Range("J85").Value = 461
$ws.Range("J86").Value = 60
$ws.Range("J88").Value = 103
$ws.Range("J90").Value = 110
$ws.Range("J91").Value = 118
$ws.Range("J94").Value = 89
$ws.Range("J95").Value = 157
$ws.Range("J99").Value = 143
$ws.Range("B101").Value = 23307
$ws.Range("D101").Value = 28142
$ws.Range("I101").Value = 26204
$ws.Range("J101").Value = 10023

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 126
$ws.Range("J7").Value = 414

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 169
$ws.Range("J4").Value = 32
$ws.Range("J7").Value = 576

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 90
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 317

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 175
$ws.Range("J7").Value = 461

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 86
$ws.Range("J6").Value = 200
$ws.Range("J7").Value = 392

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 32
$ws.Range("J3").Value = 34
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J2").Value = 4
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 84
$ws.Range("J3").Value = 110
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 67
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I4").Value = 6
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J3").Value = 18
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 196
$ws.Range("J3").Value = 204
$ws.Range("J4").Value = 34
$ws.Range("J6").Value = 183
$ws.Range("J7").Value = 635

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 100
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 301
